# ---------------------------------------------------------------------------
# Appends 10 new match rows (rows 212-221, i.e. data records 211-220 / Rodada 22)
# to the "Brazil Serie A_2022" results table on Sheet1, extending the used range
# from A1:AO211 to A1:AO221. Existing rows 1-211 are left untouched.
#
# Each element of $newRows is one spreadsheet row, itself an array of 41 cell
# descriptors @(kind, value) for columns A..AO, where kind is "n" (number) or
# "s" (text).
# ---------------------------------------------------------------------------

$newRows = @(
    @(@("n","211"),@("s","Brazil Serie A"),@("s","2022"),@("n","44786.6875"),@("n","22"),@("s","Goiás"),@("s","Avaí"),@("n","2.7"),@("n","1.95"),@("n","4.4"),@("n","1.5"),@("n","3.25"),@("n","9.25"),@("n","2.4"),@("n","1.29"),@("n","1.05"),@("n","0"),@("n","1"),@("n","1"),@("n","2.03"),@("n","3.2"),@("n","3.78"),@("n","1.08"),@("n","1.44"),@("n","2.21"),@("n","7"),@("n","2.62"),@("n","1.59"),@("n","1"),@("n","1"),@("n","2"),@("n","2"),@("n","1.71"),@("s","['79']"),@("s","['45']"),@("n","5"),@("n","4"),@("n","15"),@("n","9"),@("n","20"),@("n","13")),
    @(@("n","212"),@("s","Brazil Serie A"),@("s","2022"),@("n","44786.79166666666"),@("n","22"),@("s","Corinthians"),@("s","Palmeiras"),@("n","3.5"),@("n","1.91"),@("n","3.2"),@("n","1.55"),@("n","3.4"),@("n","9.4"),@("n","2.39"),@("n","1.3"),@("n","1.05"),@("n","0"),@("n","0"),@("n","0"),@("n","3.18"),@("n","3.05"),@("n","2.34"),@("n","1.09"),@("n","1.44"),@("n","2.27"),@("n","8.15"),@("n","2.82"),@("n","1.56"),@("n","0"),@("n","1"),@("n","1"),@("n","1.99"),@("n","1.78"),@("s","[]"),@("s","['72']"),@("n","4"),@("n","5"),@("n","12"),@("n","6"),@("n","16"),@("n","11")),
    @(@("n","213"),@("s","Brazil Serie A"),@("s","2022"),@("n","44786.85416666666"),@("n","22"),@("s","Cuiabá"),@("s","Juventude"),@("n","2.88"),@("n","1.91"),@("n","3.9"),@("n","1.61"),@("n","3.7"),@("n","10"),@("n","2.26"),@("n","1.26"),@("n","1.04"),@("n","1"),@("n","0"),@("n","1"),@("n","2.15"),@("n","2.95"),@("n","3.76"),@("n","1.11"),@("n","1.5"),@("n","2.69"),@("n","7.3"),@("n","2.55"),@("n","1.41"),@("n","1"),@("n","0"),@("n","1"),@("n","2.17"),@("n","1.66"),@("s","['22']"),@("s","[]"),@("n","5"),@("n","3"),@("n","8"),@("n","3"),@("n","13"),@("n","6")),
    @(@("n","214"),@("s","Brazil Serie A"),@("s","2022"),@("n","44786.875"),@("n","22"),@("s","Botafogo"),@("s","Atlético GO"),@("n","2.93"),@("n","2.09"),@("n","4.44"),@("n","1.51"),@("n","3.1"),@("n","8.25"),@("n","2.62"),@("n","1.33"),@("n","1.06"),@("n","0"),@("n","0"),@("n","0"),@("n","2.01"),@("n","3.12"),@("n","3.37"),@("n","1.08"),@("n","1.42"),@("n","2.18"),@("n","8.5"),@("n","2.85"),@("n","1.73"),@("n","0"),@("n","0"),@("n","0"),@("n","1.9"),@("n","1.8"),@("s","[]"),@("s","[]"),@("n","6"),@("n","4"),@("n","11"),@("n","10"),@("n","17"),@("n","14")),
    @(@("n","215"),@("s","Brazil Serie A"),@("s","2022"),@("n","44787.45833333334"),@("n","22"),@("s","Coritiba"),@("s","Atlético Mineiro"),@("n","4.5"),@("n","1.98"),@("n","2.55"),@("n","1.47"),@("n","3.2"),@("n","8.75"),@("n","2.45"),@("n","1.3"),@("n","1.06"),@("n","0"),@("n","0"),@("n","0"),@("n","4.1"),@("n","3.15"),@("n","1.85"),@("n","1.08"),@("n","1.4"),@("n","2.1"),@("n","7"),@("n","2.75"),@("n","1.65"),@("n","0"),@("n","1"),@("n","1"),@("n","2"),@("n","1.71"),@("s","[]"),@("s","['90+4']"),@("n","3"),@("n","10"),@("n","5"),@("n","12"),@("n","8"),@("n","22")),
    @(@("n","216"),@("s","Brazil Serie A"),@("s","2022"),@("n","44787.66666666666"),@("n","22"),@("s","Flamengo"),@("s","Atlético PR"),@("n","1.95"),@("n","2.3"),@("n","6.25"),@("n","1.36"),@("n","2.65"),@("n","6.75"),@("n","2.87"),@("n","1.41"),@("n","1.09"),@("n","0"),@("n","0"),@("n","0"),@("n","1.41"),@("n","4.2"),@("n","6.6"),@("n","1.05"),@("n","1.29"),@("n","1.75"),@("n","9"),@("n","3.5"),@("n","1.95"),@("n","5"),@("n","0"),@("n","5"),@("n","2.05"),@("n","1.7"),@("s","['56', '59', '63', '71', '90+2']"),@("s","[]"),@("n","12"),@("n","4"),@("n","9"),@("n","2"),@("n","21"),@("n","6")),
    @(@("n","217"),@("s","Brazil Serie A"),@("s","2022"),@("n","44787.66666666666"),@("n","22"),@("s","São Paulo"),@("s","Bragantino"),@("n","2.7"),@("n","2"),@("n","4.1"),@("n","1.46"),@("n","3"),@("n","8.25"),@("n","2.55"),@("n","1.33"),@("n","1.06"),@("n","1"),@("n","0"),@("n","1"),@("n","1.96"),@("n","3.2"),@("n","3.65"),@("n","1.08"),@("n","1.4"),@("n","2.1"),@("n","7"),@("n","2.75"),@("n","1.65"),@("n","3"),@("n","0"),@("n","3"),@("n","1.93"),@("n","1.78"),@("s","['25', '59', '61']"),@("s","[]"),@("n","6"),@("n","5"),@("n","12"),@("n","9"),@("n","18"),@("n","14")),
    @(@("n","218"),@("s","Brazil Serie A"),@("s","2022"),@("n","44787.66666666666"),@("n","22"),@("s","Ceará"),@("s","Fortaleza"),@("n","3.1"),@("n","1.93"),@("n","3.7"),@("n","1.51"),@("n","3.3"),@("n","9.5"),@("n","2.35"),@("n","1.28"),@("n","1.05"),@("n","0"),@("n","1"),@("n","1"),@("n","2.4"),@("n","2.95"),@("n","2.9"),@("n","1.11"),@("n","1.5"),@("n","2.35"),@("n","6"),@("n","2.5"),@("n","1.53"),@("n","0"),@("n","1"),@("n","1"),@("n","2.05"),@("n","1.7"),@("s","[]"),@("s","['17']"),@("n","5"),@("n","6"),@("n","9"),@("n","7"),@("n","14"),@("n","13")),
    @(@("n","219"),@("s","Brazil Serie A"),@("s","2022"),@("n","44787.75"),@("n","22"),@("s","América Mineiro"),@("s","Santos"),@("n","2.95"),@("n","1.93"),@("n","3.9"),@("n","1.51"),@("n","3.3"),@("n","9.25"),@("n","2.37"),@("n","1.29"),@("n","1.05"),@("n","1"),@("n","0"),@("n","1"),@("n","2.27"),@("n","3.1"),@("n","3"),@("n","1.1"),@("n","1.44"),@("n","2.3"),@("n","6.5"),@("n","2.62"),@("n","1.55"),@("n","1"),@("n","0"),@("n","1"),@("n","2"),@("n","1.73"),@("s","['14']"),@("s","[]"),@("n","6"),@("n","2"),@("n","12"),@("n","11"),@("n","18"),@("n","13")),
    @(@("n","220"),@("s","Brazil Serie A"),@("s","2022"),@("n","44787.79166666666"),@("n","22"),@("s","Internacional"),@("s","Fluminense"),@("n","3.15"),@("n","1.95"),@("n","3.5"),@("n","1.5"),@("n","3.25"),@("n","10"),@("n","2.45"),@("n","1.3"),@("n","1.05"),@("n","1"),@("n","0"),@("n","1"),@("n","2.49"),@("n","3"),@("n","2.75"),@("n","1.09"),@("n","1.42"),@("n","2.25"),@("n","8.35"),@("n","2.87"),@("n","1.57"),@("n","3"),@("n","0"),@("n","3"),@("n","1.95"),@("n","1.78"),@("s","['36', '71', '90+3']"),@("s","[]"),@("n","7"),@("n","0"),@("n","10"),@("n","10"),@("n","17"),@("n","10"))
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 212      # first row to insert into (row 212 == data record 211)
$lastDataRow = 211      # last pre-existing data row, used as the style template
$colCount    = 41       # columns A (1) .. AO (41)
$seasonCol   = 3        # column C ("Season") holds text that looks numeric ("2022")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstNewRow + $i
    $rowData = $newRows[$i]

    for ($c = 1; $c -le $colCount; $c++) {
        $kind = $rowData[$c - 1][0]
        $val  = $rowData[$c - 1][1]
        $cell = $ws.Cells.Item($r, $c)

        if ($kind -eq "s") {
            if ($c -eq $seasonCol) {
                # "2022" etc. would otherwise be auto-coerced to a number by Excel;
                # format as Text first, then reset the format so no visible style sticks.
                $cell.NumberFormat = "@"
                $cell.Value = $val
                $cell.Style = "Normal"
            } else {
                $cell.Value = $val
            }
        } else {
            $cell.Value = [double]$val
        }
    }

    # Match the existing table's per-column formatting: column A carries the bold/
    # bordered "Nº" style, column D carries the date/time number format. Copy those
    # formats (not values) from the last pre-existing row so new rows look identical.
    $ws.Cells.Item($lastDataRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($lastDataRow, 4).Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

Write-Output "Inserted $($newRows.Count) rows ($firstNewRow..$($firstNewRow + $newRows.Count - 1)) into '$($ws.Name)'; used range is now $($ws.UsedRange.Address())"
